$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_Checkout")

$ws.Range("L2").Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range("L3").Value = 'VerifyElement: null'
$ws.Range("L4").Value = 'Click: null'
$ws.Range("L5").Value = 'VerifyElement: null'
$ws.Range("L6").Value = 'SetText: Randomemailid'
$ws.Range("L7").Value = 'SetText: 123456'
$ws.Range("L8").Value = 'Click: null'
$ws.Range("L9").Value = 'VerifyText: Akash sangal'
$ws.Range("L10").Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range("L11").Value = 'VerifyElement: null'
$ws.Range("L12").Value = 'MoveToProductList: Quick & Easy Food Solutions'
$ws.Range("L13").Value = 'MoveAndAddProduct: null'
$ws.Range("M13").Value = 'Pass'
$ws.Range("L14").Value = 'MoveAndAddProduct: exist'
$ws.Range("M14").Value = 'Pass'
$ws.Range("N14").Value = '-'
$ws.Range("L15").Value = 'MoveAndAddProduct: Banquet Brown ''N Serve Turkey Sausage Links'
$ws.Range("L16").Value = 'MoveAndAddProduct: exist'
$ws.Range("L17").Value = 'MoveAndAddProduct: null'
$ws.Range("L18").Value = 'MoveAndAddProduct: exist'
$ws.Range("L19").Value = 'MoveAndAddProduct: College Inn Fat Free & Lower Sodium Chicken Broth - 32oz'
$ws.Range("L20").Value = 'MoveAndAddProduct: exist'
$ws.Range("L21").Value = 'MoveAndAddProduct: College Inn Fat Free & Lower Sodium Chicken Broth - 32oz'
$ws.Range("L22").Value = 'MoveAndAddProduct: exist'
$ws.Range("L23").Value = 'MoveAndAddProduct: null'
$ws.Range("M23").Value = 'Fail'
$ws.Range("L24").Value = 'Click: null'
$ws.Range("L25").Value = 'VerifyElement: null'
$ws.Range("L26").Value = 'Click: null'
$ws.Range("L27").Value = 'SelectSubstitute: Do Not Allow Substitute'
$ws.Range("L28").Value = 'SelectSubstitute: Allow Substitute'
$ws.Range("L29").Value = 'SelectSubstitute: Allow Substitute'
$ws.Range("L30").Value = 'SelectSubstitute: Do Not Allow Substitute'
$ws.Range("L31").Value = 'SelectSubstitute: Allow Substitute'
$ws.Range("L32").Value = 'SelectSubstitute: Allow Substitute'
$ws.Range("L33").Value = 'VerifySummaryViewCart: null'
$ws.Range("L34").Value = 'Wait: 5000'
$ws.Range("L35").Value = 'Click: null'
$ws.Range("L36").Value = 'VerifyTitle: Checkout'
$ws.Range("L37").Value = 'Wait: 6000'
$ws.Range("L38").Value = 'Click: null'
$ws.Range("L39").Value = 'VerifyText: Please select delivery slot.'
$ws.Range("L40").Value = 'SetText: ChandnaUserAddress'
$ws.Range("L41").Value = 'SetText: ChandnaUserCity'
$ws.Range("L42").Value = 'SetText: 068202322'
$ws.Range("L43").Value = 'SetText: 1234567890'
$ws.Range("N44").Value = 'no such element: Unable to locate element: {"method":"xpath","selector":"//textarea[normalize-space(@placeholder) = ''Order Instructions'']"}
  (Session info: chrome=68.0.3440.106)
  (Driver info: chromedriver=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91),platform=Windows NT 10.0.17134 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 38 milliseconds
For documentation on this error, please visit: http://seleniumhq.org/exceptions/no_such_element.html
Build info: version: ''unknown'', revision: ''1969d75'', time: ''2016-10-18 09:43:45 -0700''
System info: host: ''DESKTOP-OEL817D'', ip: ''192.168.134.2'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''1.8.0_161''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91), userDataDir=C:\Users\Akash\AppData\Local\Temp\scoped_dir13564_28830}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=68.0.3440.106, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=false, acceptInsecureCerts=false, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true, setWindowRect=true, unexpectedAlertBehaviour=}]
Session ID: 3fddb53966380820a83efadd593af767
*** Element info: {Using=xpath, value=//textarea[normalize-space(@placeholder) = ''Order Instructions'']}'
$ws.Range("L45").Value = 'Wait: 6000'
$ws.Range("N46").Value = 'no such element: Unable to locate element: {"method":"xpath","selector":"//table[normalize-space(@class) = ''table time-slot'']"}
  (Session info: chrome=68.0.3440.106)
  (Driver info: chromedriver=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91),platform=Windows NT 10.0.17134 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 64 milliseconds
For documentation on this error, please visit: http://seleniumhq.org/exceptions/no_such_element.html
Build info: version: ''unknown'', revision: ''1969d75'', time: ''2016-10-18 09:43:45 -0700''
System info: host: ''DESKTOP-OEL817D'', ip: ''192.168.134.2'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''1.8.0_161''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91), userDataDir=C:\Users\Akash\AppData\Local\Temp\scoped_dir13564_28830}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=68.0.3440.106, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=false, acceptInsecureCerts=false, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true, setWindowRect=true, unexpectedAlertBehaviour=}]
Session ID: 3fddb53966380820a83efadd593af767
*** Element info: {Using=xpath, value=//table[normalize-space(@class) = ''table time-slot'']}'
$ws.Range("L47").Value = 'SetText: 06820'
$ws.Range("L48").Value = 'Wait: 6000'
$ws.Range("L49").Value = 'VerifyElement: null'
$ws.Range("L50").Value = 'SelectSlot: 4 Sep'
$ws.Range("M50").Value = 'Fail'
$ws.Range("L51").Value = 'Click: null'
$ws.Range("L52").Value = 'action click failed'
$ws.Range("M52").Value = 'Fail'
$ws.Range("N52").Value = 'element not visible
  (Session info: chrome=68.0.3440.106)
  (Driver info: chromedriver=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91),platform=Windows NT 10.0.17134 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 37 milliseconds
Build info: version: ''unknown'', revision: ''1969d75'', time: ''2016-10-18 09:43:45 -0700''
System info: host: ''DESKTOP-OEL817D'', ip: ''192.168.134.2'', os.name: ''Windows 10'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''1.8.0_161''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.36.540470 (e522d04694c7ebea4ba8821272dbef4f9b818c91), userDataDir=C:\Users\Akash\AppData\Local\Temp\scoped_dir13564_28830}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=68.0.3440.106, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=false, acceptInsecureCerts=false, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true, setWindowRect=true, unexpectedAlertBehaviour=}]
Session ID: 3fddb53966380820a83efadd593af767'
$ws.Range("L53").Value = 'VerifyText: Your card''s security code is incomplete'
$ws.Range("M53").Value = 'Fail'
$ws.Range("L54").Value = 'Click: null'
$ws.Range("L55").Value = 'DeleteProductfromViewCart: all'
$ws.Range("L56").Value = 'Click: null'
$ws.Range("L57").Value = 'Click: null'
$ws.Range("L58").Value = 'Wait: 6000'
$ws.Range("L59").Value = 'VerifyElement: null'

# Writing the long multi-paragraph error texts above triggers this runtime's
# row-autofit, which the source diff does not exhibit (row heights stay at
# their original 41.25). Restore the affected rows' heights explicitly.
$ws.Rows(44).RowHeight = 41.25
$ws.Rows(46).RowHeight = 41.25
$ws.Rows(52).RowHeight = 41.25
